$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New task rows appended below the existing list (rows 6-15, gap, rows 17-21)
$ws.Range("A6").Value  = "Complete basic migration setup"
$ws.Range("A7").Value  = "Adjust career model file with pk and relationships"
$ws.Range("A8").Value  = "Create career service"
$ws.Range("A9").Value  = "Add in dummy methods for reading careers, creating careers and editing careers"
$ws.Range("A10").Value = "Add logic for creating a career"
$ws.Range("A11").Value = "Test"
$ws.Range("A12").Value = "Debug"
$ws.Range("A13").Value = "Add logic for reading careers"
$ws.Range("A14").Value = "Test"
$ws.Range("A15").Value = "Debug"

$ws.Range("A17").Value = "Create a user service"
$ws.Range("A18").Value = "Add methods shells for basic crud operations"
$ws.Range("A19").Value = "Add create method logic"
$ws.Range("A20").Value = "Test"
$ws.Range("A21").Value = "Debug"

# Mark the completed setup/migration tasks (A3:A9) with a green highlight,
# matching the "Green, Accent 6, Lighter 60%" theme swatch used in Excel.
$done = $ws.Range("A3:A9")
$done.Interior.ThemeColor = 10
$done.Interior.TintAndShade = 0.6

# The career-service tasks (A10:A15) had their fill explicitly cleared back
# to "No Fill" after being touched, rather than never having been formatted.
$cleared = $ws.Range("A10:A15")
$cleared.Interior.ColorIndex = -4142

# Move the active selection to the next empty row, ready for the next entry.
$ws.Range("A22").Select()
